$wb = $excel.ActiveWorkbook
$records = $wb.Worksheets.Item("RECORDS")
$logs = $wb.Worksheets.Item("LOGS")

# Fix the status of row 8 (spam) to Inactive (0)
$records.Range("N8").Value = 0

# Append new log entries to the LOGS sheet
$newLogs = @(
    @("ram", "Successfully Logged In!", 45792, 45792.137465277781),
    @("ram", "Marked 'asdsgtrbtynbyd' as Inactive.", 45792, 45792.13758101852),
    @("ram", "Marked 'asdasghjkl;' as Inactive.", 45792, 45792.13758101852),
    @("ram", "Reactivated user 'asdsgtrbtynbyd'.", 45792, 45792.137881944444),
    @("ram", "Reactivated user 'asdasghjkl;'.", 45792, 45792.137881944444),
    @("ram", "Reactivated user 'new'.", 45792, 45792.137881944444),
    @("ram", "Marked 'asdsgtrbtynbyd' as Inactive.", 45792, 45792.138009259259),
    @("ram", "Successfully Logged In!", 45792, 45792.138877314814),
    @("ram", "Marked 'asdasghjkl;' as Inactive.", 45792, 45792.138993055552),
    @("ram", "Marked 'spam' as Inactive.", 45792, 45792.139155092591),
    @("ram", "Marked 'new' as Inactive.", 45792, 45792.139155092591),
    @("ram", "Reactivated user 'asdsgtrbtynbyd'.", 45792, 45792.139270833337),
    @("ram", "Reactivated user 'asdasghjkl;'.", 45792, 45792.139270833337)
)

$startRow = 151
for ($i = 0; $i -lt $newLogs.Count; $i++) {
    $r = $startRow + $i
    $row = $newLogs[$i]
    $logs.Cells.Item($r, 1).Value = $row[0]
    $logs.Cells.Item($r, 2).Value = $row[1]
    $logs.Cells.Item($r, 3).Value = $row[2]
    $logs.Cells.Item($r, 3).NumberFormat = "M/d/yyyy"
    $logs.Cells.Item($r, 4).Value = $row[3]
    $logs.Cells.Item($r, 4).NumberFormat = "M/d/yyyy"
}
